# Actualización desde MV -datos-
# Append the two new daily "Dólar observado" records (04-10-2021 and 05-10-2021)
# to the bottom of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 190 -> 191).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row1 = $lastRow + 1
$row2 = $lastRow + 2

# The date column stores plain text like "01-10-2021", not real dates.
# Pre-format the cells as Text so Excel doesn't auto-convert the strings
# into date serial numbers, then clear that formatting again afterwards
# so the new cells end up with the workbook's default (unstyled) look,
# matching every other data row.
$dateRange = $ws.Range("A" + $row1 + ":A" + $row2)
$dateRange.NumberFormat = "@"

$ws.Range("A" + $row1).Value = "04-10-2021"
$ws.Range("B" + $row1).Value = 803.9

$ws.Range("A" + $row2).Value = "05-10-2021"
$ws.Range("B" + $row2).Value = 805.89

$dateRange.ClearFormats()
